# Update der excel tabelle
# Add a new experiment-log row to "Tabelle1" on the "Versuchsprotokoll (MOLE)" sheet,
# documenting the first real test-bench run, and tidy up the "Kommentar" column
# formatting (header + long comment cells get word-wrap / top alignment).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Versuchsprotokoll (MOLE)")
$lo = $ws.ListObjects.Item("Tabelle1")

# Grow the table by one row (this also extends ref/autoFilter to A1:L4).
$lo.ListRows.Add() | Out-Null

# Seed the new row from the row above so number formats / date formatting
# carry over, then overwrite the actual values.
$ws.Range("A3:L3").Copy($ws.Range("A4:L4"))

$ws.Range("A4").Value2 = 46034
$ws.Range("B4").Value = "Run_01_serial"
$ws.Range("C4").Value = "Serial (feedforward)"
$ws.Range("D4").Value = "Trajectory_02"
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = 3
$ws.Range("G4").Value = "relative_2"
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = "Stochastic"
$ws.Range("J4").Value = "Simulation"
$ws.Range("K4").Value = "nein"
$ws.Range("L4").Value = "Erster Testlauf am Prüfstand, Funktion zur Nullung des Eingangssignals vor Versuchsbeginn in Simulink erst mittendrin hinzugefügt"

# Kommentar column: wrap the long text and align it to the top, same as the
# header gets a wrap-capable style.
$ws.Range("L1").WrapText = $true

$ws.Range("L3").WrapText = $true
$ws.Range("L3").VerticalAlignment = -4160

$ws.Range("L4").WrapText = $true
$ws.Range("L4").VerticalAlignment = -4160

# The new comment is long, give row 4 enough height to show it wrapped.
$ws.Range("A4").EntireRow.RowHeight = 72

$ws.Range("K14").Select() | Out-Null

$wb.Save()
